$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.886.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.138.59'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.128.53'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.499'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.93'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +14.02%  '
$ws.Range('E11').Value = '  +2.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.464'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.92'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.09%  '
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.642.44'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.937.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '527.71'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +10.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.138.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.44%  '
$ws.Range('E20').Value = '  +3.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.702'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.41'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '78.49'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +14.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.80'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.41%  '
$ws.Range('E29').Value = '  +2.81%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.12'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('E33').Value = '  +4.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '561.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.38'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0440'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '52.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0815'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.79%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.066.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.22%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +15.70%  '
$ws.Range('E42').Value = '  +3.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('E44').Value = '  +7.25%  '
$ws.Range('E45').Value = '  +8.53%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '25.02'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '119.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.109'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.52%  '
$ws.Range('E51').Value = '  +4.07%  '
